# Junction_Flooding_460.xlsx refresh: "custom accuracy + 데이터 1000개"
# - Replace the 4 data rows (rows 2-5) with the newly-sampled readings.
# - Drop the now-unused 5th data row (old row 6) from the sheet.
# - A handful of numeric columns grew one character wider (7 -> 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset now only has 4 sample rows, so remove the old trailing row 6
# (this also shifts the used range / dimension down to A1:AH5).
$ws.Rows.Item(6).Delete()

# Newly sampled Time + J1..J33 readings for rows 2-5.
$rows = @(
    @(45141.50694444445, 18.256, 12.153, 4.015, 38.832, 30.967, 14.367, 45.024, 22.106, 9.170999999999999, 13.751, 15.281, 15.885, 4.586, 14.287, 19.896, 12.285, 3.423, 2.238, 210.098, 39.697, 13.187, 26.04, 13.224, 3.103, 22.935, 11.648, 10.595, 12.441, 15.796, 3.454, 39.912, 7.199, 16.486),
    @(45141.51388888889, 15.854, 11.191, 1.792, 34.293, 27.658, 12.477, 48.12, 19.197, 8.286, 12.186, 13.714, 14.332, 3.986, 12.407, 17.476, 10.708, 1.538, 1.058, 181.527, 34.726, 11.452, 23.006, 11.934, 2.255, 23.435, 10.115, 9.153, 10.735, 14.26, 1.265, 43.798, 6.315, 14.318),
    @(45141.52083333334, 14.893, 10.723, 1.271, 32.312, 26.205, 11.721, 46.04, 18.033, 7.883, 11.587, 12.959, 13.59, 3.744, 11.655, 16.484, 9.99, 1.031, 0.779, 170.08, 32.615, 10.758, 21.724, 11.337, 1.952, 22.205, 9.502000000000001, 8.535, 10.016, 13.527, 0.773, 41.786, 5.972, 13.45),
    @(45141.52777777778, 13.93, 10.12, 1.02, 30.26, 24.61, 10.96, 43.23, 16.87, 7.42, 10.91, 12.14, 12.76, 3.5, 10.9, 15.45, 9.300000000000001, 0.79, 0.65, 158.63, 30.51, 10.06, 20.38, 10.66, 1.75, 20.84, 8.890000000000001, 7.95, 9.34, 12.7, 0.5600000000000001, 39.23, 5.61, 12.58)
)

$numRows = $rows.Count
$numCols = $rows[0].Count
$data = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $rows[$r][$c]
    }
}

$ws.Range("A2:AH5").Value = $data

# Columns C, G, Q, V, AA, AB, AC widened from 7 to 8 characters.
# (ColumnWidth is expressed in characters; 7.1667 round-trips to a stored
# raw column width of 8 for this workbook's default font.)
$wideColumns = @(3, 7, 17, 22, 27, 28, 29)
foreach ($col in $wideColumns) {
    $ws.Columns.Item($col).ColumnWidth = 7.1667
}
